# hits_map_name.xlsx - fix numbering mistake
# The "Rank"/numbering column (A) in both tables on sheet "Plate1" was
# numbered 1..12 top-to-bottom; it should instead count down 12..1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1: rows 4-15, column A ---
$table1Rows = 4..15
$n = $table1Rows.Count
for ($i = 0; $i -lt $n; $i++) {
    $row = $table1Rows[$i]
    $ws.Range("A$row").Value = $n - $i
}

# --- Table 2: rows 20-31, column A ---
$table2Rows = 20..31
$n2 = $table2Rows.Count
for ($i = 0; $i -lt $n2; $i++) {
    $row = $table2Rows[$i]
    $ws.Range("A$row").Value = $n2 - $i
}

# --- Update the selection/view: move to J40 (and this also clears the
#     previous topLeftCell scroll-freeze at A8) ---
$ws.Activate() | Out-Null
$ws.Range("J40").Select() | Out-Null
